$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "release version" column (V) values from "V1" to "V1.0" for the
#    existing reviewed rows.
$ws.Range("V2").Value = "V1.0"
$ws.Range("V3").Value = "V1.0"
$ws.Range("V4").Value = "V1.0"

# 2. Add a new row (row 5) capturing the customer's response for the
#    SIQ_user_transactions requirements.
$ws.Range("O5").Value = "SIQ_user_transactions_06`nSIQ_user_transactions_07`nSIQ_user_transactions_08`ntheses SRS should be deleted as we can mention all the bank account details in SIQ_user_transactions_01"
$ws.Range("P5").Value = "Aya"
$ws.Range("Q5").Value = "Sara"
$ws.Range("R5").Value = "User Transactions"
$ws.Range("S5").Value = "V3.0"
$ws.Range("T5").Value = "Sara"
$ws.Range("U5").Value = "User Transactions"
$ws.Range("V5").Value = "V3.0"

# Match formatting used by the other rows: the comment cell (O) is
# center aligned with wrapped text, the remaining cells are simply
# center aligned.
$ws.Range("O5").HorizontalAlignment = -4108
$ws.Range("O5").VerticalAlignment = -4108
$ws.Range("O5").WrapText = $true

$ws.Range("P5:V5").HorizontalAlignment = -4108
$ws.Range("P5:V5").VerticalAlignment = -4108

# Size the new row to fit the wrapped comment text.
$ws.Rows.Item(5).RowHeight = 87

# Reflect the reviewer's last active cell selection.
$ws.Range("V2").Select()
